# Add two new columns (I: "I0", J: "IF") to the data table, mirroring the
# formatting of the existing header cells and filling in the data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting (bold font, border, center/top alignment) from the
# existing "IP" header cell (H1) into the two new header cells, then set
# their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-15) ---
$iValues = @(7, 7, 6, 8, 9, 7, 7, 8, 7, 6, 5, 9, 9, 8)
$jValues = @(8, 9, 6, 9, 9, 8, 8, 9, 8, 6, 6, 9, 9, 8)

for ($r = 2; $r -le 15; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
